# Medidas.xlsx - add new animation-frame measurement blocks (Saltar avanzando,
# Saltar estatico, Puñetazo simple, Puñetazo fuerte, Gancho) below the
# existing Idle / Caminar tables, matching rows 15-37.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 15-23: "Saltar avanzando" block --------------------------------
$ws.Range("A15").Value = "Saltar avanzando"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = 32
$ws.Range("D15").Value = 186
$ws.Range("E15").Formula = "=85-C15"
$ws.Range("F15").Formula = "=268-D15"

$ws.Range("B16").Formula = "=1+B15"
$ws.Range("C16").Value = 125
$ws.Range("D16").Value = 156
$ws.Range("E16").Formula = "=166-C16"
$ws.Range("F16").Formula = "=263-D16"

$ws.Range("B17").Formula = "=1+B16"
$ws.Range("C17").Value = 199
$ws.Range("D17").Value = 156
$ws.Range("E17").Formula = "=253-C17"
$ws.Range("F17").Formula = "=232-D17"

$ws.Range("B18").Formula = "=1+B17"
$ws.Range("C18").Value = 283
$ws.Range("D18").Value = 159
$ws.Range("E18").Formula = "=352-C18"
$ws.Range("F18").Formula = "=244-D18"

$ws.Range("B19").Formula = "=1+B18"
$ws.Range("C19").Value = 377
$ws.Range("D19").Value = 171
$ws.Range("E19").Formula = "=458-C19"
$ws.Range("F19").Formula = "=223-D19"

$ws.Range("B20").Formula = "=1+B19"
$ws.Range("C20").Value = 488
$ws.Range("D20").Value = 161
$ws.Range("E20").Formula = "=552-C20"
$ws.Range("F20").Formula = "=221-D20"

$ws.Range("B21").Formula = "=1+B20"
$ws.Range("C21").Value = 582
$ws.Range("D21").Value = 167
$ws.Range("E21").Formula = "=663-C21"
$ws.Range("F21").Formula = "=219-D21"

$ws.Range("B22").Formula = "=1+B21"
$ws.Range("C22").Value = 697
$ws.Range("D22").Value = 167
$ws.Range("E22").Formula = "=761-C22"
$ws.Range("F22").Formula = "=227-D22"

$ws.Range("B23").Formula = "=1+B22"
$ws.Range("C23").Value = 32
$ws.Range("D23").Value = 186
$ws.Range("E23").Formula = "=85-C23"
$ws.Range("F23").Formula = "=268-D23"

# --- Row 24-28: "Saltar estatico" block ----------------------------------
$ws.Range("A24").Value = "Saltar estatico"
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = 32
$ws.Range("D24").Value = 186
$ws.Range("E24").Formula = "=85-C24"
$ws.Range("F24").Formula = "=268-D24"

$ws.Range("B25").Formula = "=1+B24"
$ws.Range("C25").Value = 125
$ws.Range("D25").Value = 156
$ws.Range("E25").Formula = "=166-C25"
$ws.Range("F25").Formula = "=263-D25"

$ws.Range("B26").Formula = "=1+B25"
$ws.Range("C26").Value = 199
$ws.Range("D26").Value = 156
$ws.Range("E26").Formula = "=253-C26"
$ws.Range("F26").Formula = "=232-D26"

$ws.Range("B27").Formula = "=1+B26"
$ws.Range("C27").Value = 283
$ws.Range("D27").Value = 159
$ws.Range("E27").Formula = "=352-C27"
$ws.Range("F27").Formula = "=244-D27"

$ws.Range("B28").Formula = "=1+B27"
$ws.Range("C28").Value = 32
$ws.Range("D28").Value = 186
$ws.Range("E28").Formula = "=85-C28"
$ws.Range("F28").Formula = "=268-D28"

# Rows 24 & 25 carry an explicit (slightly taller) custom row height.
$ws.Rows.Item(24).RowHeight = 15.75
$ws.Rows.Item(25).RowHeight = 15.75

# --- Row 29-30: "Puñetazo simple" block ----------------------------------
$ws.Range("A29").Value = "Puñetazo simple"
$ws.Range("B29").Value = 1
$ws.Range("C29").Value = 29
$ws.Range("D29").Value = 316
$ws.Range("E29").Formula = "=92-C29"
$ws.Range("F29").Formula = "=405-D29"

$ws.Range("B30").Value = 2
$ws.Range("C30").Value = 127
$ws.Range("D30").Value = 313
$ws.Range("E30").Formula = "=210-C30"
$ws.Range("F30").Formula = "=402-D30"

# --- Row 31-32: "Puñetazo fuerte" block ----------------------------------
$ws.Range("A31").Value = "Puñetazo fuerte"
$ws.Range("B31").Value = 1
$ws.Range("C31").Value = 237
$ws.Range("D31").Value = 321
$ws.Range("E31").Formula = "=304-C31"
$ws.Range("F31").Formula = "=401-D31"

$ws.Range("B32").Value = 2
$ws.Range("C32").Value = 333
$ws.Range("D32").Value = 324
$ws.Range("E32").Formula = "=442-C32"
$ws.Range("F32").Formula = "=401-D32"

# --- Row 33-37: "Gancho" block -------------------------------------------
$ws.Range("A33").Value = "Gancho"
$ws.Range("B33").Value = 1
$ws.Range("C33").Value = 465
$ws.Range("D33").Value = 320
$ws.Range("E33").Formula = "=532-C33"
$ws.Range("F33").Formula = "=400-D33"

$ws.Range("B34").Value = 2
$ws.Range("C34").Value = 559
$ws.Range("D34").Value = 334
$ws.Range("E34").Formula = "=637-C34"
$ws.Range("F34").Formula = "=403-D34"

$ws.Range("B35").Value = 3
$ws.Range("C35").Value = 668
$ws.Range("D35").Value = 325
$ws.Range("E35").Formula = "=750-C35"
$ws.Range("F35").Formula = "=402-D35"

$ws.Range("B36").Value = 4
$ws.Range("C36").Value = 774
$ws.Range("D36").Value = 315
$ws.Range("E36").Formula = "=853-C36"
$ws.Range("F36").Formula = "=405-D36"

$ws.Range("B37").Value = 5
$ws.Range("C37").Value = 874
$ws.Range("D37").Value = 293
$ws.Range("E37").Formula = "=929-C37"
$ws.Range("F37").Formula = "=403-D37"

# Give the new (taller, multi-row) animation headers the same bold
# center/center look as the existing ones, but with word-wrap turned on
# so the longer names fit the column. Build the wrap-enabled format once
# on a scratch cell, then stamp it onto every header before merging (so
# every cell in each merged block carries the right style index).
$ws.Range("A3").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("Z1").WrapText = $true

$ws.Range("Z1").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A31").PasteSpecial(-4122)
$ws.Range("A33").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("Z1").Clear()

# --- Merge the animation-name column for each new block ------------------
$ws.Range("A15:A23").Merge()
$ws.Range("A24:A28").Merge()
$ws.Range("A29:A30").Merge()
$ws.Range("A31:A32").Merge()
$ws.Range("A33:A37").Merge()

# --- Reposition the view / selection so the freshly added rows are
#     visible, matching where the author ended up while entering data.
$ws.Range("F28").Select()
